# Journal de Travail - add two new journal entries (rows 46 and 47)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33 already carries the alternate-week style pattern (A/B=s1, C-G=s2, H=s8)
# that rows 46 & 47 need, so copy its formatting down first.
$ws.Range("A33:H33").Copy()
$ws.Range("A46:H46").PasteSpecial(-4122)
$ws.Range("A33:H33").Copy()
$ws.Range("A47:H47").PasteSpecial(-4122)

$ws.Rows.Item(46).RowHeight = 47.25
$ws.Rows.Item(47).RowHeight = 35.25

# Row 46: Documentation / Mise a jour du rapport
$ws.Range("A46").Value = 44266
$ws.Range("B46").Value = "Documentationt"
$ws.Range("C46").Value = "6h"
$ws.Range("D46").Value = "Mise a jour du rapport"
$ws.Range("E46").Value = "Oui"
$ws.Range("F46").Value = "Permet d'attribuer des annonces au joueurs pendant une partie"
$ws.Range("G46").Value = "Non"
$ws.Range("H46").Value = ""

# Row 47: Rails (Backend) / Mise en place du serveur heroku
$ws.Range("A47").Value = 44266
$ws.Range("B47").Value = "Rails (Backend)"
$ws.Range("C47").Value = "2h"
$ws.Range("D47").Value = "Mise en place du serveur heroku"
$ws.Range("E47").Value = "Oui"
$ws.Range("F47").Value = "Permet d'avoir accès a l'API n'importe ou"
$ws.Range("G47").Value = "Non"
$ws.Range("H47").Value = ""

$ws.Range("A47").Select()
